$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the "Periodo Mora" (period) column E for rows 16-25 to the new
# ascending order (2008 .. 2105) instead of the old descending order.
$periods = @("2008", "2009", "2010", "2011", "2012", "2101", "2102", "2103", "2104", "2105")
for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = 16 + $i
    $ws.Range("E$row").Value = $periods[$i]
}

# Update "Salario Basico" (column G) for the whole table: every row moves
# from 828116 to 877803.
for ($row = 16; $row -le 25; $row++) {
    $ws.Range("G$row").Value = 877803
}

# Update "Valor Mora" (column F). All rows keep 35112 except the row that
# now corresponds to period 2105 (row 25), which keeps the special value
# 28090 that used to belong to period 2105 back when it sat on row 16.
for ($row = 16; $row -le 24; $row++) {
    $ws.Range("F$row").Value = 35112
}
$ws.Range("F25").Value = 28090
